# Launch service provider notes - proofing-pass style edit.
#
# Applies (in document order):
#   1. Nanoracks heading       -> wrap run in spellStart/spellEnd proofErr
#   2. "51.6 degree inclination" -> split into two runs w/ gramStart/gramEnd
#      around "51.6 degree"
#   3. SpaceX paragraph        -> split run so "Cubesats" is wrapped in
#      spellStart/spellEnd
#   4. GOMspace heading        -> wrap run in spellStart/spellEnd proofErr
#   5. GOMspace description    -> split run so "5 year" is wrapped in
#      gramStart/gramEnd
#   6. Endurosat heading       -> wrap run in spellStart/spellEnd proofErr
#   7. New paragraph containing the NASA CubeSat diagram URL, inserted
#      right after the picture paragraph (before "1U = 10cm x 10cm x 10cm")
#   8. "Typically ~ 1kg-1.3kg per 1U" -> split so "Typically" is wrapped in
#      gramStart/gramEnd
#
# Each edit replaces a whole paragraph's Range via Range.InsertXML with the
# exact WordprocessingML desired, which lets us plant <w:proofErr/> markers
# that aren't reachable through the higher-level object-model properties.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# ---------------------------------------------------------------------
# 1. Nanoracks
# ---------------------------------------------------------------------
Set-ParagraphXml 1 @'
<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Nanoracks</w:t></w:r><w:proofErr w:type="spellEnd"/>
'@

# ---------------------------------------------------------------------
# 2. 51.6 degree inclination
# ---------------------------------------------------------------------
Set-ParagraphXml 5 @'
<w:proofErr w:type="gramStart"/><w:r><w:t>51.6 degree</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> inclination</w:t></w:r>
'@

# ---------------------------------------------------------------------
# 3. SpaceX rideshare paragraph
# ---------------------------------------------------------------------
Set-ParagraphXml 8 @'
<w:r><w:t>Minimum 300</w:t></w:r><w:r><w:t xml:space="preserve">k$ for a rideshare, might be the best option for launching multiple </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cubesats</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as other providers might charge per satellite rather than by a minimum + extra for mass like SpaceX do.</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; We may need to package them together and plan for separation.</w:t></w:r>
'@

# ---------------------------------------------------------------------
# 4. GOMspace heading
# ---------------------------------------------------------------------
Set-ParagraphXml 9 @'
<w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>GOMspace</w:t></w:r><w:proofErr w:type="spellEnd"/>
'@

# ---------------------------------------------------------------------
# 5. GOMspace description
# ---------------------------------------------------------------------
Set-ParagraphXml 10 @'
<w:r><w:t xml:space="preserve">Offers better integrated packages for </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>5 year</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> LEO missions with all required hardware &#8211; not an option for our mission.</w:t></w:r>
'@

# ---------------------------------------------------------------------
# 6. Endurosat heading
# ---------------------------------------------------------------------
Set-ParagraphXml 15 @'
<w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Endurosat</w:t></w:r><w:proofErr w:type="spellEnd"/>
'@

# ---------------------------------------------------------------------
# 7. New paragraph with the NASA CubeSat diagram URL, right after the
#    embedded picture paragraph and before "1U = 10cm x 10cm x 10cm".
# ---------------------------------------------------------------------
$picPara = $d.Paragraphs.Item(19)
$picPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(20)
$newPara.Range.InsertXML("<w:p $wNs><w:r><w:t>https://www.nasa.gov/wp-content/uploads/2015/03/what_are_cubesats.png</w:t></w:r></w:p>")

# ---------------------------------------------------------------------
# 8. Typically ~ 1kg-1.3kg per 1U.
#    (the trailing "." run is left untouched; it was already split out)
# ---------------------------------------------------------------------
Set-ParagraphXml 22 @'
<w:proofErr w:type="gramStart"/><w:r><w:t>Typically</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ~ 1kg-1.3kg per 1U</w:t></w:r><w:r><w:t>.</w:t></w:r>
'@

Write-Output "ok"
